$wb = $excel.ActiveWorkbook

# --- CaseToForm sheet (sheet2): add a mandatory "id" column (row 4: id / `this` / note) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E4").Value = "id"
$ws2.Range("F4").Value = '`this`'
$ws2.Range("G4").Value = "In order to sync with a SQL database, id is mandatory"

# New column G needs to be wide enough to show the note text (~26.29 chars, matches the
# existing bespoke widths used on columns A/C of this sheet).
$ws2.Columns.Item(7).ColumnWidth = 25.57

# CaseToForm becomes the active/selected sheet, with G5 as the new selection anchor.
[void]$ws2.Activate()
[void]$ws2.Range("G5").Select()
